$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Predicted_next_Day_Price" column (AB) values for rows 2-7
$ws.Range("AB2").Value = 0.5758531837252193
$ws.Range("AB3").Value = 0.3875584177971766
$ws.Range("AB4").Value = 0.263323795011487
$ws.Range("AB5").Value = 0.184042929938158
$ws.Range("AB6").Value = 0.3585390729305011
$ws.Range("AB7").Value = 0.5927845411628769

# Update the "Predicted_Signal" column (AC) values that shifted
$ws.Range("AC5").Value = 0
$ws.Range("AC6").Value = 1
$ws.Range("AC7").Value = 0

# Update the "Actual_Return" column (AD) values that shifted
$ws.Range("AD6").Value = 0
$ws.Range("AD7").Value = 0.001870779339295581
